$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 updates
$ws.Range("G5").Value = 2.4
$ws.Range("I5").Value = 3.1
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("AV5").Value = 4.75
$ws.Range("AY5").Value = 51

# Row 7 updates
$ws.Range("G7").Value = 1.7
$ws.Range("H7").Value = 3.55
$ws.Range("I7").Value = 4.5
$ws.Range("J7").Value = 2.27
$ws.Range("K7").Value = 2.12
$ws.Range("L7").Value = 4.8
$ws.Range("N7").Value = 9.449999999999999
$ws.Range("O7").Value = 1.34
$ws.Range("P7").Value = 2.75
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 1.65
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.5
$ws.Range("X7").Value = 7.2
$ws.Range("Z7").Value = 12.5
$ws.Range("AA7").Value = 15
$ws.Range("AB7").Value = 32
$ws.Range("AD7").Value = 7
$ws.Range("AG7").Value = 10.75
$ws.Range("AH7").Value = 24
$ws.Range("AI7").Value = 15.5
$ws.Range("AJ7").Value = 75
$ws.Range("AN7").Value = 3.4
$ws.Range("AO7").Value = 8.25
$ws.Range("AP7").Value = 19.5
$ws.Range("AQ7").Value = 29
$ws.Range("AS7").Value = 2.47
$ws.Range("AV7").Value = 6.1
$ws.Range("AW7").Value = 26
